# Atualização de bases das ligas, do dia: 10-06-2024 às 07:08
#
# The refreshed feed re-sorted a few fixtures that share the same date,
# so the rows holding their id / scores / odds need to be exchanged
# (column A = row index, C = Div, D = Date stay the same for a pair,
# since they are identical between the two rows anyway).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows {
    param($RowA, $RowB)

    # Swap column B (id)
    $idA = $ws.Range("B$RowA").Value2
    $idB = $ws.Range("B$RowB").Value2
    $ws.Range("B$RowA").Value = $idB
    $ws.Range("B$RowB").Value = $idA

    # Swap columns E:AD (HomeTeam .. PL_AhUnder)
    $rangeA = $ws.Range("E${RowA}:AD${RowA}")
    $rangeB = $ws.Range("E${RowB}:AD${RowB}")

    $valsA = $rangeA.Value2
    $valsB = $rangeB.Value2

    $rangeA.Value = $valsB
    $rangeB.Value = $valsA
}

Swap-Rows 124 125
Swap-Rows 168 169
Swap-Rows 194 195
